$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Checksum" field row (row 15 of the "Documents[].File" group) was removed.
# Deleting the entire row shifts all following rows up by one and Excel
# automatically adjusts the merged cell ranges in columns A and B.
$ws.Rows.Item(15).Delete()
